# Update the "Metadata" worksheet (sheet 1) of the ValueSet workbook:
#  - Version 0.1.6 -> 0.1.7
#  - Status active -> draft
#  - Date updated
#  - Contact row gets the full publisher contact string
#  - A new "Contact" row (Bob Milius) + a new "Jurisdiction" row are inserted
#  - everything below shifts down by one row
#  - the "Include from RxNorm" worksheet (sheet 2) needs no content changes

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Make room for the new "Jurisdiction" row ------------------------------
# Shift rows 12..15 (Description, Purpose, Copyright, Immutable) down to
# rows 13..16, working from the bottom up so we never overwrite a value
# before it has been copied.
$ws1.Range("A16").Value = $ws1.Range("A15").Value()
$ws1.Range("B16").Value = $ws1.Range("B15").Value()

$ws1.Range("A15").Value = $ws1.Range("A14").Value()
$ws1.Range("B15").Value = $ws1.Range("B14").Value()

$ws1.Range("A14").Value = $ws1.Range("A13").Value()
$ws1.Range("B14").Value = $ws1.Range("B13").Value()

$ws1.Range("A13").Value = $ws1.Range("A12").Value()
$ws1.Range("B13").Value = $ws1.Range("B12").Value()

# Copy the row formatting (border/fill/font/alignment) down onto row 16,
# which previously did not exist.
$ws1.Range("A15:B15").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 10: first Contact row gets the real contact string ----------------
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Row 11: second Contact row -> the new Jurisdiction placeholder --------
$ws1.Range("A11").Value = "Contact"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Row 12: brand new "Jurisdiction" row (empty value) ---------------------
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""

# --- Simple value updates ----------------------------------------------------
$ws1.Range("B3").Value = "0.1.7"
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2024-08-27T12:23:18-05:00"

Write-Output "Metadata sheet updated"
